$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header text updates (shared-string runs collapse to plain text on write,
# but font/size/color are identical across runs so rendering is unaffected)
$ws.Range("A8").Value = "Volume 30   Number  16"
$ws.Range("C9").Value = "Report Covering the Week  4/17/2023  Through  4/23/2023"

# Weekly crime statistics table updates (rows 14-30, columns C-N)
# Row 14
$ws.Range("C14").Value = 6
$ws.Range("D14").Value = 4
$ws.Range("E14").Value = 50
$ws.Range("F14").Value = 32
$ws.Range("G14").Value = 28
$ws.Range("H14").Value = 14.285714285714
$ws.Range("I14").Value = 120
$ws.Range("J14").Value = 125
$ws.Range("K14").Value = -4
$ws.Range("L14").Value = -11.111111111111
$ws.Range("M14").Value = -17.241379310344
$ws.Range("N14").Value = -79.487179487179

# Row 15
$ws.Range("C15").Value = 29
$ws.Range("D15").Value = 27
$ws.Range("E15").Value = 7.407407407407
$ws.Range("F15").Value = 113
$ws.Range("G15").Value = 115
$ws.Range("H15").Value = -1.739130434782
$ws.Range("I15").Value = 461
$ws.Range("J15").Value = 493
$ws.Range("K15").Value = -6.490872210953
$ws.Range("L15").Value = 7.209302325581
$ws.Range("M15").Value = 24.258760107816
$ws.Range("N15").Value = -51.061571125265

# Row 16
$ws.Range("C16").Value = 266
$ws.Range("D16").Value = 282
$ws.Range("E16").Value = -5.673758865248
$ws.Range("F16").Value = 1132
$ws.Range("G16").Value = 1198
$ws.Range("H16").Value = -5.509181969949
$ws.Range("I16").Value = 4693
$ws.Range("J16").Value = 4817
$ws.Range("K16").Value = -2.574216317209
$ws.Range("L16").Value = 43.736600306278
$ws.Range("M16").Value = -12.947505101094
$ws.Range("N16").Value = -81.821351100093

# Row 17
$ws.Range("C17").Value = 524
$ws.Range("D17").Value = 416
$ws.Range("E17").Value = 25.961538461538
$ws.Range("F17").Value = 2111
$ws.Range("G17").Value = 1949
$ws.Range("H17").Value = 8.31195484864
$ws.Range("I17").Value = 7884
$ws.Range("J17").Value = 7199
$ws.Range("K17").Value = 9.515210445895
$ws.Range("L17").Value = 34.127254168084
$ws.Range("M17").Value = 62.255608149825
$ws.Range("N17").Value = -30.118773267151

# Row 18
$ws.Range("C18").Value = 277
$ws.Range("D18").Value = 273
$ws.Range("E18").Value = 1.465201465201
$ws.Range("F18").Value = 1046
$ws.Range("G18").Value = 1173
$ws.Range("H18").Value = -10.82693947144
$ws.Range("I18").Value = 4406
$ws.Range("J18").Value = 4742
$ws.Range("K18").Value = -7.085617882749
$ws.Range("L18").Value = 24.322799097065
$ws.Range("M18").Value = -17.706387747478
$ws.Range("N18").Value = -85.49369505811

# Row 19
$ws.Range("C19").Value = 994
$ws.Range("D19").Value = 927
$ws.Range("E19").Value = 7.22761596548
$ws.Range("F19").Value = 3768
$ws.Range("G19").Value = 3644
$ws.Range("H19").Value = 3.402854006586
$ws.Range("I19").Value = 14771
$ws.Range("J19").Value = 14962
$ws.Range("K19").Value = -1.276567303836
$ws.Range("L19").Value = 53.656506813689
$ws.Range("M19").Value = 38.111266947171
$ws.Range("N19").Value = -39.118786579836

# Row 20
$ws.Range("C20").Value = 307
$ws.Range("D20").Value = 210
$ws.Range("E20").Value = 46.190476190476
$ws.Range("F20").Value = 1180
$ws.Range("G20").Value = 865
$ws.Range("H20").Value = 36.416184971098
$ws.Range("I20").Value = 4492
$ws.Range("J20").Value = 3966
$ws.Range("K20").Value = 13.262733232476
$ws.Range("L20").Value = 90.50042408821
$ws.Range("M20").Value = 51.245791245791
$ws.Range("N20").Value = -87.209567198177

# Row 21
$ws.Range("C21").Value = 2403
$ws.Range("D21").Value = 2139
$ws.Range("E21").Value = 12.342215988779
$ws.Range("F21").Value = 9382
$ws.Range("G21").Value = 8972
$ws.Range("H21").Value = 4.569772625947
$ws.Range("I21").Value = 36827
$ws.Range("J21").Value = 36304
$ws.Range("K21").Value = 1.440612604671
$ws.Range("L21").Value = 46.005629782341
$ws.Range("M21").Value = 23.64277320799
$ws.Range("N21").Value = -71.314067611777

# Row 22
$ws.Range("C22").Value = 43
$ws.Range("D22").Value = 43
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 179
$ws.Range("G22").Value = 165
$ws.Range("H22").Value = 8.484848484848
$ws.Range("I22").Value = 663
$ws.Range("J22").Value = 710
$ws.Range("K22").Value = -6.619718309859
$ws.Range("L22").Value = 56.367924528301
$ws.Range("M22").Value = 1.84331797235

# Row 23
$ws.Range("C23").Value = 123
$ws.Range("D23").Value = 118
$ws.Range("E23").Value = 4.237288135593
$ws.Range("F23").Value = 500
$ws.Range("G23").Value = 427
$ws.Range("H23").Value = 17.096018735363
$ws.Range("I23").Value = 1872
$ws.Range("J23").Value = 1701
$ws.Range("K23").Value = 10.05291005291
$ws.Range("L23").Value = 22.67365661861
$ws.Range("M23").Value = 69.10569105691

# Row 24
$ws.Range("C24").Value = 2120
$ws.Range("D24").Value = 2094
$ws.Range("E24").Value = 1.24164278892
$ws.Range("F24").Value = 8127
$ws.Range("G24").Value = 8673
$ws.Range("H24").Value = -6.295399515738
$ws.Range("I24").Value = 32601
$ws.Range("J24").Value = 32603
$ws.Range("K24").Value = -0.006134404809
$ws.Range("L24").Value = 41.798964812317
$ws.Range("M24").Value = 42.698940733607

# Row 25
$ws.Range("C25").Value = 827
$ws.Range("D25").Value = 741
$ws.Range("E25").Value = 11.605937921727
$ws.Range("F25").Value = 3323
$ws.Range("G25").Value = 3150
$ws.Range("H25").Value = 5.492063492063
$ws.Range("I25").Value = 12607
$ws.Range("J25").Value = 11947
$ws.Range("K25").Value = 5.524399430819
$ws.Range("L25").Value = 36.958174904943
$ws.Range("M25").Value = -4.143856447688

# Row 26
$ws.Range("C26").Value = 44
$ws.Range("D26").Value = 41
$ws.Range("E26").Value = 7.317073170731
$ws.Range("F26").Value = 176
$ws.Range("H26").Value = -3.825136612021
$ws.Range("I26").Value = 735
$ws.Range("J26").Value = 795
$ws.Range("K26").Value = -7.54716981132
$ws.Range("L26").Value = 1.80055401662

# Row 27
$ws.Range("C27").Value = 125
$ws.Range("D27").Value = 78
$ws.Range("E27").Value = 60.25641025641
$ws.Range("F27").Value = 421
$ws.Range("G27").Value = 373
$ws.Range("H27").Value = 12.868632707774
$ws.Range("I27").Value = 1542
$ws.Range("J27").Value = 1430
$ws.Range("K27").Value = 7.832167832167
$ws.Range("L27").Value = 20.280811232449

# Row 28
$ws.Range("C28").Value = 13
$ws.Range("D28").Value = 26
$ws.Range("E28").Value = -50
$ws.Range("F28").Value = 82
$ws.Range("G28").Value = 126
$ws.Range("H28").Value = -34.920634920634
$ws.Range("I28").Value = 336
$ws.Range("J28").Value = 436
$ws.Range("K28").Value = -22.935779816513
$ws.Range("L28").Value = -14.066496163682
$ws.Range("M28").Value = -24.8322147651
$ws.Range("N28").Value = -79.844031193761

# Row 29
$ws.Range("C29").Value = 12
$ws.Range("D29").Value = 21
$ws.Range("E29").Value = -42.857142857142
$ws.Range("F29").Value = 70
$ws.Range("G29").Value = 96
$ws.Range("H29").Value = -27.083333333333
$ws.Range("I29").Value = 281
$ws.Range("J29").Value = 372
$ws.Range("K29").Value = -24.462365591397
$ws.Range("L29").Value = -20.621468926553
$ws.Range("M29").Value = -23.224043715847
$ws.Range("N29").Value = -81.525312294543

# Row 30
$ws.Range("D30").Value = 8
$ws.Range("E30").Value = -50
$ws.Range("F30").Value = 36
$ws.Range("G30").Value = 31
$ws.Range("H30").Value = 16.129032258064
$ws.Range("I30").Value = 136
$ws.Range("J30").Value = 215
$ws.Range("K30").Value = -36.744186046511
$ws.Range("L30").Value = 4.615384615384

